$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 11 and 12 ("Paso 5 - El sistema invoca CU05..." and
# "Paso 6 - El sistema invoca CU06...") entirely.
$ws.Rows("11:12").Delete()

# The two "Invocando CUxx" steps (now rows 16 and 17) lose their
# "Invocando CUxx..." suffix text.
$ws.Range("B16").Value = "El sistema modifica el estado del ticket."
$ws.Range("B17").Value = "El sistema modifica el estado de la oportunidad."

$ws.Range("B10").Select()
